# Edit the Configuration file
# - Settings sheet: update Business Process Name, add three new settings rows
#   (System1_Credential, System1_URL, System1_WorkItemsURL) with hyperlinks
#   on the two URL values.
# - Constants sheet: bump MaxRetryNumber from 0 to 2.
# - Selection / active-tab bookkeeping to match what Excel leaves behind
#   after a user edits the Settings sheet.

$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# ---- Settings sheet -------------------------------------------------
# Business process name value (row 5, column B)
$wsSettings.Range("B5").Value = "REF_Generate Yearly Report Dispatcher"

# New rows 6-8 (fill order chosen to match the authored shared-string order)
$wsSettings.Range("A7").Value = "System1_URL"
$wsSettings.Range("A8").Value = "System1_WorkItemsURL"
$wsSettings.Range("A6").Value = "System1_Credential"
$wsSettings.Range("B8").Value = "https://acme-test.uipath.com/work-items/"
$wsSettings.Range("B7").Value = "https://acme-test.uipath.com/"
$wsSettings.Range("B6").Value = "ACME_Cred"

# Hyperlinks on the URL cells (added in this order so B8 -> rId1, B7 -> rId2)
$wsSettings.Hyperlinks.Add($wsSettings.Range("B8"), "https://acme-test.uipath.com/work-items/")
$wsSettings.Hyperlinks.Add($wsSettings.Range("B7"), "https://acme-test.uipath.com/")

# ---- Constants sheet --------------------------------------------------
# MaxRetryNumber changes from 0 to 2
$wsConstants.Range("B2").Value = 2

# ---- Selection / active sheet bookkeeping -----------------------------
[void]$wsConstants.Activate()
[void]$wsConstants.Range("B2").Select()

[void]$wsSettings.Activate()
[void]$wsSettings.Range("C19").Select()
